$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new trade row (row 17) with the same shape/style as the existing rows.
$ws.Cells.Item(17, 1).Value = 9440.75
$ws.Cells.Item(17, 2).Value = 9798.39
$ws.Cells.Item(17, 3).Value = 277
$ws.Cells.Item(17, 4).Value = 287.11
$ws.Cells.Item(17, 5).Value = $true
$ws.Cells.Item(17, 6).Value = 3.65

$ws.Cells.Item(17, 7).Value = [DateTime]::FromOADate(42626.545324074075)
$ws.Cells.Item(16, 7).Copy()
$ws.Cells.Item(17, 7).PasteSpecial(-4122)

$ws.Cells.Item(17, 8).Value = $false
